# Adding Solar Energy Ulmeni to the Portfolio Forecast
#
# The forecast window shifts forward by 3 days (new data pull), so every
# timestamp in column A (rows 2-97) moves from 2025-05-30/05-31 to
# 2025-06-02/06-03. In addition, the production forecast values for the
# morning ramp-up (rows 21-41, column B) are updated with new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original timestamps (Excel serial dates) for rows 2..97, in order.
$aVals = @(
    45807.01041666666, 45807.02083333334, 45807.03125,       45807.04166666666,
    45807.05208333334, 45807.0625,         45807.07291666666, 45807.08333333334,
    45807.09375,        45807.10416666666, 45807.11458333334, 45807.125,
    45807.13541666666, 45807.14583333334, 45807.15625,       45807.16666666666,
    45807.17708333334, 45807.1875,         45807.19791666666, 45807.20833333334,
    45807.21875,        45807.22916666666, 45807.23958333334, 45807.25,
    45807.26041666666, 45807.27083333334, 45807.28125,       45807.29166666666,
    45807.30208333334, 45807.3125,         45807.32291666666, 45807.33333333334,
    45807.34375,        45807.35416666666, 45807.36458333334, 45807.375,
    45807.38541666666, 45807.39583333334, 45807.40625,       45807.41666666666,
    45807.42708333334, 45807.4375,         45807.44791666666, 45807.45833333334,
    45807.46875,        45807.47916666666, 45807.48958333334, 45807.5,
    45807.51041666666, 45807.52083333334, 45807.53125,       45807.54166666666,
    45807.55208333334, 45807.5625,         45807.57291666666, 45807.58333333334,
    45807.59375,        45807.60416666666, 45807.61458333334, 45807.625,
    45807.63541666666, 45807.64583333334, 45807.65625,       45807.66666666666,
    45807.67708333334, 45807.6875,         45807.69791666666, 45807.70833333334,
    45807.71875,        45807.72916666666, 45807.73958333334, 45807.75,
    45807.76041666666, 45807.77083333334, 45807.78125,       45807.79166666666,
    45807.80208333334, 45807.8125,         45807.82291666666, 45807.83333333334,
    45807.84375,        45807.85416666666, 45807.86458333334, 45807.875,
    45807.88541666666, 45807.89583333334, 45807.90625,       45807.91666666666,
    45807.92708333334, 45807.9375,         45807.94791666666, 45807.95833333334,
    45807.96875,        45807.97916666666, 45807.98958333334, 45808
)

# Shift every timestamp forward by exactly 3 days and write it back.
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i] + 3
}

# New production values (column B) for the morning ramp-up period.
$bUpdates = @{
    21 = 5;    22 = 18;   23 = 38;   24 = 66;   25 = 109;  26 = 199;  27 = 283;
    28 = 361;  29 = 437;  30 = 616;  31 = 733;  32 = 840;  33 = 939;  34 = 1057;
    35 = 1149; 36 = 1255; 37 = 1280; 38 = 1411; 39 = 1476; 40 = 1547; 41 = 1586
}

foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}
